$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data to reflect the latest refresh.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.190.75"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.481.21"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.01"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.08"
$ws.Range("E6").Value = "  +2.17%  "
$ws.Range("E7").Value = "  -0.69%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.77"
$ws.Range("E10").Value = "  +6.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0807"
$ws.Range("E11").Value = "  -1.38%  "
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.13"
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.11"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.848.24"
$ws.Range("E15").Value = "  -1.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.497.39"
$ws.Range("E16").Value = "  +1.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.844"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.110.82"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.68"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0930"
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.74"
$ws.Range("E22").Value = "  +14.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.25"
$ws.Range("E23").Value = "  -0.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "244.97"
$ws.Range("E24").Value = "  -2.69%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.65"
$ws.Range("E27").Value = "  -2.73%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.00"
$ws.Range("E28").Value = "  +1.35%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.19"
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.75"
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.134"
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.42"
$ws.Range("E32").Value = "  -0.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.92"
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.34"
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0782"
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.96"
$ws.Range("E37").Value = "  +1.55%  "
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("E40").Value = "  -0.54%  "
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "119.14"
$ws.Range("E42").Value = "  -3.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.82"
$ws.Range("E43").Value = "  +4.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0294"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.979.38"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.00"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.99"
$ws.Range("E47").Value = "  -5.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.07"
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("E49").Value = "  -2.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.13"
$ws.Range("E50").Value = "  -5.16%  "
$ws.Range("E51").Value = "  +4.47%  "
